# sp_Blitz Check ID List - add new checks (#115 Parallelism Rocket Surgery,
# #182 Backup Compression Default Off/Memory Pressure, #183 Hardware NUMA,
# etc.) to the "Sheet2" (Check ID List) worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# Row 179: CheckID 114 - Hardware - NUMA Config (Server Info, no URL)
$ws.Cells.Item(179, 1).Value = 114
$ws.Cells.Item(179, 2).Value = 250
$ws.Cells.Item(179, 3).Value = "Server Info"
$ws.Cells.Item(179, 4).Value = "Hardware - NUMA Config"

# Row 180: CheckID 115 - Parallelism Rocket Surgery (Performance)
$ws.Cells.Item(180, 1).Value = 115
$ws.Cells.Item(180, 2).Value = 110
$ws.Cells.Item(180, 3).Value = "Performance"
$ws.Cells.Item(180, 4).Value = "Parallelism Rocket Surgery"
$ws.Cells.Item(180, 5).Value = "http://BrentOzar.com/go/makeparallel"
[void]$ws.Hyperlinks.Add($ws.Cells.Item(180, 5), "http://BrentOzar.com/go/makeparallel")

# Row 181: CheckID 116 - Backup Compression Default Off (Informational)
$ws.Cells.Item(181, 1).Value = 116
$ws.Cells.Item(181, 2).Value = 200
$ws.Cells.Item(181, 3).Value = "Informational"
$ws.Cells.Item(181, 4).Value = "Backup Compression Default Off"
$ws.Cells.Item(181, 5).Value = "http://BrentOzar.com/go/backup"
[void]$ws.Hyperlinks.Add($ws.Cells.Item(181, 5), "http://BrentOzar.com/go/backup")

# Row 182: CheckID 117 - Memory Pressure Affecting Queries (Performance)
$ws.Cells.Item(182, 1).Value = 117
$ws.Cells.Item(182, 2).Value = 100
$ws.Cells.Item(182, 3).Value = "Performance"
$ws.Cells.Item(182, 4).Value = "Memory Pressure Affecting Queries"
$ws.Cells.Item(182, 5).Value = "http://BrentOzar.com/go/grants"
[void]$ws.Hyperlinks.Add($ws.Cells.Item(182, 5), "http://BrentOzar.com/go/grants")

# Update the view: scroll the frozen pane down and select E183 as the new
# active cell in the bottom-right pane.
[void]$ws.Range("E183").Select()
